$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.618
$ws.Range("A3").Value = -21.573
$ws.Range("D3").Value = -7.768000000000001
$ws.Range("D12").Value = -7.260000000000001
$ws.Range("A14").Value = -21.559
$ws.Range("A16").Value = -21.373
$ws.Range("B18").Value = 5.494
$ws.Range("A21").Value = -20.664
$ws.Range("A23").Value = -20.673
$ws.Range("B24").Value = 6.055
$ws.Range("D24").Value = -7.839
$ws.Range("A25").Value = -21.275
$ws.Range("B25").Value = 7.053999999999999
$ws.Range("D25").Value = -8.146000000000001
$ws.Range("A26").Value = -21.481
$ws.Range("B27").Value = 6.202000000000001
$ws.Range("A29").Value = -21.202
$ws.Range("B30").Value = 5.656999999999999
$ws.Range("B31").Value = 6.56
$ws.Range("B39").Value = 7.717000000000001
$ws.Range("A40").Value = -20.404
$ws.Range("D41").Value = -7.944
$ws.Range("B42").Value = 8.33
$ws.Range("B48").Value = 5.274
$ws.Range("D50").Value = -8.436999999999998
$ws.Range("B51").Value = 6.718999999999999
$ws.Range("B52").Value = 6.407000000000001
$ws.Range("A53").Value = -22.04
$ws.Range("D53").Value = -7.704000000000001
$ws.Range("B55").Value = 4.790999999999999
$ws.Range("B56").Value = 6.448
$ws.Range("D56").Value = -7.867
$ws.Range("A57").Value = -21.303
$ws.Range("B57").Value = 6.431999999999999
$ws.Range("D57").Value = -8.232000000000001
$ws.Range("D58").Value = -8.096
$ws.Range("A59").Value = -22.43
$ws.Range("B60").Value = 5.787000000000001
$ws.Range("D61").Value = -7.972999999999999
$ws.Range("D63").Value = -7.529000000000001
$ws.Range("D64").Value = -7.743
$ws.Range("A65").Value = -21.418
$ws.Range("A69").Value = -21.531
$ws.Range("D70").Value = -7.498000000000002
$ws.Range("D72").Value = -7.715000000000001
$ws.Range("B73").Value = 6.679
$ws.Range("B74").Value = 8.995000000000001
$ws.Range("A79").Value = -21.136
$ws.Range("A83").Value = -21.351
$ws.Range("D86").Value = -8.273
$ws.Range("B89").Value = 5.898999999999999
$ws.Range("D89").Value = -6.093999999999999
$ws.Range("B90").Value = 6.045
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 5.898999999999999
$ws.Range("A93").Value = -21.186
$ws.Range("D98").Value = -8.134
$ws.Range("A100").Value = -21.419
$ws.Range("D100").Value = -7.989
$ws.Range("D102").Value = -7.933
